$wb = $excel.ActiveWorkbook

# --- sheet1 "credentials": change selection (no longer the active tab) ---
$ws1 = $wb.Worksheets.Item(1)
$null = $ws1.Range("F42").Select()

# --- sheet3 "sites" -> rename to "logins" and replace its data ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "logins"

# remove the hyperlinks and the extra (4th) data row entirely
$null = $ws3.Cells.Hyperlinks.Delete()
$null = $ws3.Rows.Item(4).Delete()

# strip the leftover hyperlink formatting from column B, then rewrite the
# sheet as a username/password table (like "credentials")
$ws3.Range("A1:B3").ClearFormats()

$ws3.Range("A1").Value = "username"
$ws3.Range("B1").Value = "password"
$ws3.Range("A2").Value = "kaiser.bobo"
$ws3.Range("B2").Value = "cutie.pie.bb2022"
$ws3.Range("A3").Value = "guligina.beauty"
$ws3.Range("B3").Value = "only1wife2019"

$ws3.Range("A1:B1").Font.Bold = $true

$ws3.Columns.Item(2).ColumnWidth = 31

$null = $ws3.Range("A4").Select()

# the Hyperlink cell style is now unused anywhere in the workbook
$null = $wb.Styles.Item("Hyperlink").Delete()

# --- add the new empty "Sheet1" right after "logins" and make it active ---
$newSheet = $wb.Worksheets.Add($null, $ws3)
$newSheet.Name = "Sheet1"
$null = $newSheet.Activate()
